$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Cheque Received Date" filter label to "...From" and add a
# matching "...To" filter (label + input cell) next to it, mirroring the
# existing "Cheque Number From"/"Cheque Number To" pair in row 5.
$ws.Range("A3").Value = "Cheque Received Date From"

# New label cell C3 - copy formatting from A3 (bold label style) then set text
$ws.Range("A3").Copy($ws.Range("C3"))
$ws.Range("C3").Value = "Cheque Received Date To"

# New input cell D3 - copy formatting from B3 (date input style), leave empty
$ws.Range("B3").Copy($ws.Range("D3"))
$ws.Range("D3").ClearContents()

# Widen column A slightly to fit the longer label text
$ws.Columns.Item(1).ColumnWidth = 25.6667

# Update the active selection to B3, as in the saved file
[void]$ws.Range("B3").Select()
